$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title strings (Jan 2017 -> Feb 2017)
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("B32").Value = "Source: Short-Term Energy Outlook, February 2017."

# Update underlying monthly history/forecast figures (rows 27-30, cols F:H)
# These feed the J:M shared formulas (=col-prevcol) and row 31 (=row27-SUM(row28:30))
$ws.Range("F27").Value = 19.584006213999999
$ws.Range("G27").Value = 19.842936774999998
$ws.Range("H27").Value = 20.172870081999999
$ws.Range("F28").Value = 9.2937664698999995
$ws.Range("G28").Value = 9.2875013014000007
$ws.Range("H28").Value = 9.3608492685000009
$ws.Range("F29").Value = 1.6040000847
$ws.Range("G29").Value = 1.6177085177999999
$ws.Range("H29").Value = 1.6141701643999999
$ws.Range("F30").Value = 3.8543217240000001
$ws.Range("G30").Value = 3.9336894548000001
$ws.Range("H30").Value = 4.0361051753000003

# Update monthly history (C) / forecast (D) series, rows 54-84
$ws.Range("C54").Value = 19.79928
$ws.Range("C55").Value = 19.712032000000001
$ws.Range("C56").Value = 20.130901000000001
$ws.Range("C59").Value = 19.654798
$ws.Range("C60").Value = 19.421787434999999
$ws.Range("D60").Value = "#N/A"
$ws.Range("C61").Value = 19.022804287
$ws.Range("D61").Value = 19.022804287
$ws.Range("D62").Value = 19.753070000000001
$ws.Range("D63").Value = 19.67202
$ws.Range("D64").Value = 19.484200000000001
$ws.Range("D65").Value = 19.579329999999999
$ws.Range("D66").Value = 20.023150000000001
$ws.Range("D67").Value = 20.245190000000001
$ws.Range("D68").Value = 20.314109999999999
$ws.Range("D69").Value = 20.091989999999999
$ws.Range("D70").Value = 19.913060000000002
$ws.Range("D71").Value = 19.927399999999999
$ws.Range("D72").Value = 20.08522
$ws.Range("D73").Value = 19.762989999999999
$ws.Range("D74").Value = 19.88646
$ws.Range("D75").Value = 19.872070000000001
$ws.Range("D76").Value = 19.851769999999998
$ws.Range("D77").Value = 19.859590000000001
$ws.Range("D78").Value = 20.33989
$ws.Range("D79").Value = 20.529769999999999
$ws.Range("D80").Value = 20.660900000000002
$ws.Range("D81").Value = 20.37133
$ws.Range("D82").Value = 20.22626
$ws.Range("D83").Value = 20.206489999999999
$ws.Range("D84").Value = 20.481719999999999
